$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheet "Tarefas": mark the last (previously last) task as "Em Andamento"
#    and stamp it with a date, like the rows right above it.
# ---------------------------------------------------------------------
$wsTarefas = $wb.Worksheets.Item("Tarefas")

# Copy the date style from A56 (same block) onto A57, then give it a value.
$wsTarefas.Range("A56").Copy()
$wsTarefas.Range("A57").PasteSpecial(-4122)
$wsTarefas.Range("A57").Value = 42204

# Status goes from "Aguardando" to "Em Andamento"
$wsTarefas.Range("B57").Value = "Em Andamento"

$wsTarefas.Activate()
$wsTarefas.Range("B58:D60").Select()

# ---------------------------------------------------------------------
# 2) Sheet "Anotações Gerais": add a new bullet point.
# ---------------------------------------------------------------------
$wsNotas = $wb.Worksheets.Item("Anotações Gerais")
$wsNotas.Range("A5").Value = "Fazer a otimização sequencial das funções?"

$wsNotas.Activate()
$wsNotas.Range("A5").Select()

# ---------------------------------------------------------------------
# 3) New sheet "mutações" with a table of JS optimization practices.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMut = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsMut.Name = "mutações"

# Cell values are entered in the same order the shared-string table records
# them (column by column, each list pasted in before its source link was
# added), so that the underlying string table comes out identical.
$wsMut.Range("A2").Value = "Retirar globais"
$wsMut.Range("A3").Value = "Retirar var de dentro de for"
$wsMut.Range("A1").Value = "http://pt.slideshare.net/doris1/performance-optimization-and-javascript-best-practices"

$wsMut.Range("B1").Value = "http://jonraasch.com/blog/10-javascript-performance-boosting-tips-from-nicholas-zakas"
$wsMut.Range("B2").Value = "Apagar instruções com with()"
$wsMut.Range("B3").Value = "Trocar for in por for puro"
$wsMut.Range("B4").Value = "Trocar for por while com variavel de controle externa (item 7)"

$wsMut.Range("C2").Value = "Use === Instead of =="
$wsMut.Range("C1").Value = "http://blogs.msdn.com/b/dorischen/archive/2011/03/17/web-performance-tips-10-javascript-best-practices.aspx"
$wsMut.Range("C3").Value = "Eval = Bad"
$wsMut.Range("C4").Value = "Reduce Globals: Namespace"
$wsMut.Range("C5").Value = "Don't Pass a String to ""SetInterval"" or ""SetTimeOut"""

$wsMut.Range("A4").Value = "Use {} Instead of New Object()"
$wsMut.Range("A5").Value = "Use [] Instead of New Array()"

$wsMut.Range("C6").Value = "Use {} Instead of New Object()"
$wsMut.Range("C7").Value = "Use [] Instead of New Array()"

# Highlight the duplicated "good practice" cells with an accent fill
# (single call over the multi-area range so only one new cell style is
# recorded instead of one per cell).
$wsMut.Range("A4,A5,C6,C7").Interior.ThemeColor = 7

$wsMut.Columns.Item(1).ColumnWidth = 81.42578125
$wsMut.Columns.Item(2).ColumnWidth = 81.85546875
$wsMut.Columns.Item(3).ColumnWidth = 104.5703125

$wsMut.Activate()
$wsMut.Range("C6").Select()
